$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 35 (old rows 35-46 shift down to 37-48)
$ws.Range("A35:R36").Insert()

# Populate the two newly inserted rows (35 and 36) with the new weekly data
$ws.Range("A35").Value = 9
$ws.Range("B35").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C35").Value = "Metropolitana"
$ws.Range("D35").Value = 44466
$ws.Range("E35").Value = 13
$ws.Range("F35").Value = 100114002
$ws.Range("G35").Value = "Camote"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 1330
$ws.Range("K35").Value = 10000
$ws.Range("L35").Value = 11000
$ws.Range("M35").Value = 10500
$ws.Range("N35").Value = "$/malla 18 kilos"
$ws.Range("O35").Value = "Perú"
$ws.Range("P35").Value = 583
$ws.Range("Q35").Value = 18
$ws.Range("R35").Value = "Hortaliza"

$ws.Range("A36").Value = 9
$ws.Range("B36").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C36").Value = "Metropolitana"
$ws.Range("D36").Value = 44466
$ws.Range("E36").Value = 13
$ws.Range("F36").Value = 100114002
$ws.Range("G36").Value = "Camote"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Segunda"
$ws.Range("J36").Value = 790
$ws.Range("K36").Value = 9000
$ws.Range("L36").Value = 9000
$ws.Range("M36").Value = 9000
$ws.Range("N36").Value = "$/malla 18 kilos"
$ws.Range("O36").Value = "Perú"
$ws.Range("P36").Value = 500
$ws.Range("Q36").Value = 18
$ws.Range("R36").Value = "Hortaliza"
